$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 8500
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 8500
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 8500
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -10746

$ws.Range("H89").Value = 8500
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 8500
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 42500
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -53732

$ws.Range("H138").Value = 3094.25
$ws.Range("I138").Value = 1676.8889
$ws.Range("J138").Value = 3944.6667
$ws.Range("K138").Value = 5030.6667
$ws.Range("L138").Value = 11834.0001
$ws.Range("M138").Value = 109.3333000000002
$ws.Range("N138").Value = -22114.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5172.129
$ws.Range("I32").Value = 5172.129
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 5172.129
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4885.129

$ws.Range("H61").Value = 2849.1428
$ws.Range("I61").Value = 2490.6667
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2490.6667
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2278.6667

$ws.Range("H74").Value = 7580.6875
$ws.Range("I74").Value = 8008.7
$ws.Range("J74").Value = 5440.625
$ws.Range("K74").Value = 8008.7
$ws.Range("L74").Value = 5440.625
$ws.Range("M74").Value = -7134.7
$ws.Range("N74").Value = -7188.625

$ws.Range("H77").Value = 7580.6875
$ws.Range("I77").Value = 8008.7
$ws.Range("J77").Value = 5440.625
$ws.Range("K77").Value = 40043.5
$ws.Range("L77").Value = 27203.125
$ws.Range("M77").Value = -35675.5
$ws.Range("N77").Value = -35939.125

$ws.Range("H102").Value = 863.2
$ws.Range("I102").Value = 933
$ws.Range("J102").Value = 409.5
$ws.Range("K102").Value = 933
$ws.Range("L102").Value = 409.5
$ws.Range("M102").Value = 689
$ws.Range("N102").Value = -3653.5

$ws.Range("H122").Value = 3781.2856
$ws.Range("I122").Value = 3828.1667
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 11484.5001
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -9034.500100000001

$ws.Range("H136").Value = 2849.1428
$ws.Range("I136").Value = 2490.6667
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 7472.000100000001
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -4922.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4309.1665
$ws.Range("I20").Value = 1171.2
$ws.Range("J20").Value = 19999
$ws.Range("K20").Value = 1171.2
$ws.Range("L20").Value = 19999
$ws.Range("M20").Value = -924.2
$ws.Range("N20").Value = -20493

$ws.Range("H105").Value = 2999.75
$ws.Range("I105").Value = 4000
$ws.Range("J105").Value = 1999.5
$ws.Range("K105").Value = 4000
$ws.Range("L105").Value = 1999.5
$ws.Range("M105").Value = -2253
$ws.Range("N105").Value = -5493.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()

$ws.Range("H58").Value = 4734.75
$ws.Range("I58").Value = 4756.1816
$ws.Range("J58").Value = 4499
$ws.Range("K58").Value = 4756.1816
$ws.Range("L58").Value = 4499
$ws.Range("M58").Value = -4553.1816
$ws.Range("N58").Value = -4905

$ws.Range("H62").Value = 4753
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4753
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 4753
$ws.Range("N62").Value = -6001

$ws.Range("H65").Value = 4753
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4753
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 23765
$ws.Range("N65").Value = -30005

$ws.Range("H125").Value = 20000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 20000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 20000
$ws.Range("N125").Value = -24920

$ws.Range("H132").Value = 3559.4546
$ws.Range("I132").Value = 3044.875
$ws.Range("J132").Value = 4931.6665
$ws.Range("K132").Value = 9134.625
$ws.Range("L132").Value = 14794.9995
$ws.Range("M132").Value = -6604.625

$ws.Range("H134").Value = 10375.857
$ws.Range("I134").Value = 10668.25
$ws.Range("J134").Value = 9986
$ws.Range("K134").Value = 32004.75
$ws.Range("L134").Value = 29958
$ws.Range("M134").Value = -29469.75

$ws.Range("H136").Value = 4734.75
$ws.Range("I136").Value = 4756.1816
$ws.Range("J136").Value = 4499
$ws.Range("K136").Value = 14268.5448
$ws.Range("L136").Value = 13497
$ws.Range("M136").Value = -11718.5448
$ws.Range("N136").Value = -18597

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 182.625
$ws.Range("I12").Value = 203
$ws.Range("J12").Value = 166.77777
$ws.Range("K12").Value = 609
$ws.Range("L12").Value = 500.33331
$ws.Range("M12").Value = -436
$ws.Range("N12").Value = -846.33331

$ws.Range("H36").Value = 187.25
$ws.Range("I36").Value = 74.5
$ws.Range("J36").Value = 300
$ws.Range("K36").Value = 223.5
$ws.Range("L36").Value = 900
$ws.Range("M36").Value = -54.5
$ws.Range("N36").Value = -1238

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4752.75
$ws.Range("I102").Value = 4752.75
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4752.75
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -3130.75
$ws.Range("N102").ClearContents()

$ws.Range("H132").Value = 3175.0667
$ws.Range("I132").Value = 2884.4546
$ws.Range("J132").Value = 3974.25
$ws.Range("K132").Value = 8653.363799999999
$ws.Range("L132").Value = 11922.75
$ws.Range("M132").Value = -6123.363799999999
$ws.Range("N132").Value = -16982.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3233.8823
$ws.Range("I132").Value = 2564.4443
$ws.Range("J132").Value = 3987
$ws.Range("K132").Value = 7693.3329
$ws.Range("L132").Value = 11961
$ws.Range("M132").Value = -5163.3329
$ws.Range("N132").Value = -17021

$ws.Range("H136").Value = 2717.5454
$ws.Range("I136").Value = 2699.3
$ws.Range("J136").Value = 2900
$ws.Range("K136").Value = 8097.900000000001
$ws.Range("L136").Value = 8700
$ws.Range("M136").Value = -5547.900000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 812.5
$ws.Range("I107").Value = 700
$ws.Range("J107").Value = 880
$ws.Range("K107").Value = 2100
$ws.Range("L107").Value = 2640
$ws.Range("M107").Value = -180
$ws.Range("N107").Value = -6480

$ws.Range("H122").Value = 3489.111
$ws.Range("I122").Value = 3489.111
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10467.333
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8017.332999999999

$ws.Range("H132").Value = 2056.7
$ws.Range("I132").Value = 852.8570999999999
$ws.Range("J132").Value = 4865.6665
$ws.Range("K132").Value = 2558.5713
$ws.Range("L132").Value = 14596.9995
$ws.Range("M132").Value = -28.57129999999961
$ws.Range("N132").Value = -19656.9995

$ws.Range("H136").Value = 3378.2
$ws.Range("I136").Value = 2766.5264
$ws.Range("J136").Value = 15000
$ws.Range("K136").Value = 8299.5792
$ws.Range("L136").Value = 45000
$ws.Range("M136").Value = -5749.5792
